# "implementing x scrolling + decoy correction"
#
# The "Block type" / "block letter" table (columns A:B, rows 4-27) is
# re-sorted alphabetically by block letter (a-z). Three new letters that
# didn't previously have a row (j, n, q) get their own placeholder rows
# (letter only, no block type yet assigned - "decoys" waiting to be
# filled in). The long-standing duplicate letter bug is also fixed:
# "hidden area center" incorrectly shared the letter "y" with "decoy";
# it is now correctly moved to use the unused letter "z", while "decoy"
# keeps "y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the existing (unsorted) A:B table body before rewriting it in the
# new sorted order - the table grows from 24 rows (4-27) to 26 (4-29).
$ws.Range("A4:B27").ClearContents()

$ws.Range("A4").Value = "Armor potion"
$ws.Range("B4").Value = "a"
$ws.Range("A5").Value = "woodBridge"
$ws.Range("B5").Value = "b"
$ws.Range("A6").Value = "Cave entrance"
$ws.Range("B6").Value = "c"
$ws.Range("A7").Value = "Door"
$ws.Range("B7").Value = "d"
$ws.Range("A8").Value = "cave Exit"
$ws.Range("B8").Value = "e"
$ws.Range("A9").Value = "Arena centering target"
$ws.Range("B9").Value = "f"
$ws.Range("A10").Value = "flaG"
$ws.Range("B10").Value = "g"
$ws.Range("A11").Value = "Health potion"
$ws.Range("B11").Value = "h"
$ws.Range("A12").Value = "trIgger"
$ws.Range("B12").Value = "i"
$ws.Range("B13").Value = "j"
$ws.Range("A14").Value = "piKes"
$ws.Range("B14").Value = "k"
$ws.Range("A15").Value = "Lava puddle"
$ws.Range("B15").Value = "l"
$ws.Range("A16").Value = "hidden area trigger"
$ws.Range("B16").Value = "m"
$ws.Range("B17").Value = "n"
$ws.Range("A18").Value = "water flOw"
$ws.Range("B18").Value = "o"
$ws.Range("A19").Value = "Pizza box"
$ws.Range("B19").Value = "p"
$ws.Range("B20").Value = "q"
$ws.Range("A21").Value = "cave Rock"
$ws.Range("B21").Value = "r"
$ws.Range("A22").Value = "Stalagmite"
$ws.Range("B22").Value = "s"
$ws.Range("A23").Value = "Tree"
$ws.Range("B23").Value = "t"
$ws.Range("A24").Value = "water pUddle"
$ws.Range("B24").Value = "u"
$ws.Range("A25").Value = "rock"
$ws.Range("B25").Value = "v"
$ws.Range("A26").Value = "Wall"
$ws.Range("B26").Value = "w"
$ws.Range("A27").Value = "Exit to next map"
$ws.Range("B27").Value = "x"
$ws.Range("A28").Value = "decoy"
$ws.Range("B28").Value = "y"
$ws.Range("A29").Value = "hidden area center"
$ws.Range("B29").Value = "z"

# Scroll/selection ended up back at the top-left (no more topLeftCell
# offset) with the active cell on the new last row of the table.
$ws.Range("A1").Select()
$ws.Range("B29").Select()
